# ------------------------------------------------------------------
# Refresh the cryptocurrency snapshot table (cryptos list) in place:
#  - columns D (Price) and E (Volume 1h) are refreshed for every coin
#  - rows 45-47 additionally rotate which coin (name/link in B/C)
#    occupies that rank, per the latest source data
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force literal text so numeric-looking strings (e.g. "1.00",
    # "0.600") keep their exact formatting instead of being
    # auto-coerced into a Number by Excel, then drop the helper
    # number format so no stray cell style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "61.545.60"
$ws.Cells.Item(2, 5).Value = "  +2.00%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "2.655.44"
$ws.Cells.Item(3, 5).Value = "  +1.92%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.03%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "581.08"
$ws.Cells.Item(5, 5).Value = "  -1.25%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "145.30"
$ws.Cells.Item(6, 5).Value = "  +1.60%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.998"
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.600"
$ws.Cells.Item(8, 5).Value = "  +0.29%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "6.58"
$ws.Cells.Item(9, 5).Value = "  +0.92%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "0.109"
$ws.Cells.Item(10, 5).Value = "  +3.44%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.380"
$ws.Cells.Item(11, 5).Value = "  +2.83%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "0.157"
$ws.Cells.Item(12, 5).Value = "  +0.98%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "3.126.95"
$ws.Cells.Item(13, 5).Value = "  +1.93%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "26.58"
$ws.Cells.Item(14, 5).Value = "  +7.86%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "61.508.48"
$ws.Cells.Item(15, 5).Value = "  +1.94%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +3.47%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "2.667.50"
$ws.Cells.Item(17, 5).Value = "  +2.16%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "11.63"
$ws.Cells.Item(18, 5).Value = "  +2.25%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "4.77"
$ws.Cells.Item(19, 5).Value = "  +2.45%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "354.98"
$ws.Cells.Item(20, 5).Value = "  +2.12%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "6.93"
$ws.Cells.Item(21, 5).Value = "  +0.34%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "0.998"
$ws.Cells.Item(22, 5).Value = "  -0.18%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "0.525"
$ws.Cells.Item(23, 5).Value = "  +0.36%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "64.15"
$ws.Cells.Item(24, 5).Value = "  +1.70%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "8.60"
$ws.Cells.Item(25, 5).Value = "  +6.90%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +2.82%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "0.998"
$ws.Cells.Item(27, 5).Value = "  -0.05%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +7.19%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "0.0₃0822"
$ws.Cells.Item(29, 5).Value = "  +3.49%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "6.86"
$ws.Cells.Item(30, 5).Value = "  +7.73%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.07%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "166.81"
$ws.Cells.Item(32, 5).Value = "  +1.99%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "20.09"
$ws.Cells.Item(33, 5).Value = "  +2.93%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "4.72"
$ws.Cells.Item(34, 5).Value = "  +10.53%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "1.10"
$ws.Cells.Item(35, 5).Value = "  +12.43%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "1.34"
$ws.Cells.Item(36, 5).Value = "  +8.39%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +6.89%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "345.85"
$ws.Cells.Item(38, 5).Value = "  +11.47%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +6.13%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.915"
$ws.Cells.Item(40, 5).Value = "  +8.94%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "38.34"
$ws.Cells.Item(41, 5).Value = "  +1.05%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "5.34"
$ws.Cells.Item(42, 5).Value = "  +6.61%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "0.0580"
$ws.Cells.Item(43, 5).Value = "  +5.48%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "21.23"
$ws.Cells.Item(44, 5).Value = "  +5.26%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Mantle"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Cells.Item(45, 4) "0.630"
$ws.Cells.Item(45, 5).Value = "  +4.26%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(46, 4) "20.56"
$ws.Cells.Item(46, 5).Value = "  +3.96%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(47, 4) "135.16"
$ws.Cells.Item(47, 5).Value = "  -0.41%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +4.47%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "0.0999"
$ws.Cells.Item(49, 5).Value = "  +0.50%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "0.998"
$ws.Cells.Item(50, 5).Value = "  +0.16%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "2.105.18"
$ws.Cells.Item(51, 5).Value = "  +3.57%  "

